$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Qty"
$ws.Range("D1").Value = "Price"
$ws.Range("E1").Value = "Total"
